$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 450.73077
$ws.Range("I8").Value = 1050.5
$ws.Range("J8").Value = 307.92856
$ws.Range("K8").Value = 3151.5
$ws.Range("L8").Value = 923.78568
$ws.Range("M8").Value = -3012.5
$ws.Range("N8").Value = -1201.78568
$ws.Range("H64").Value = 9501.5
$ws.Range("J64").Value = 9003
$ws.Range("L64").Value = 9003
$ws.Range("N64").Value = -9499
$ws.Range("H67").Value = 9501.5
$ws.Range("J67").Value = 9003
$ws.Range("L67").Value = 9003
$ws.Range("N67").Value = -10719
$ws.Range("H74").Value = 8077.8
$ws.Range("I74").Value = 6445.75
$ws.Range("J74").Value = 9165.833000000001
$ws.Range("K74").Value = 6445.75
$ws.Range("L74").Value = 9165.833000000001
$ws.Range("M74").Value = -5509.75
$ws.Range("N74").Value = -11037.833
$ws.Range("H77").Value = 8077.8
$ws.Range("I77").Value = 6445.75
$ws.Range("J77").Value = 9165.833000000001
$ws.Range("K77").Value = 32228.75
$ws.Range("L77").Value = 45829.165
$ws.Range("M77").Value = -27548.75
$ws.Range("N77").Value = -55189.165
$ws.Range("I111").Value = 31258752
$ws.Range("K111").Value = 93776256
$ws.Range("M111").Value = -93773189
$ws.Range("H113").Value = 33341056
$ws.Range("J113").Value = 62511740
$ws.Range("L113").Value = 62511740
$ws.Range("N113").Value = -62518248
$ws.Range("H116").Value = 15629488
$ws.Range("I116").Value = 31252000
$ws.Range("J116").Value = 6975.625
$ws.Range("K116").Value = 31252000
$ws.Range("L116").Value = 6975.625
$ws.Range("M116").Value = -31248558
$ws.Range("N116").Value = -13859.625
$ws.Range("H132").Value = 1320.8334
$ws.Range("I132").Value = 1516.238
$ws.Range("K132").Value = 4548.714
$ws.Range("M132").Value = -2018.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 9999.75
$ws.Range("J3").Value = 9999.75
$ws.Range("L3").Value = 9999.75
$ws.Range("N3").Value = -10229.75
$ws.Range("H63").Value = 1797.5
$ws.Range("I63").Value = 1797.5
$ws.Range("K63").Value = 1797.5
$ws.Range("M63").Value = -1111.5
$ws.Range("H66").Value = 1797.5
$ws.Range("I66").Value = 1797.5
$ws.Range("K66").Value = 8987.5
$ws.Range("M66").Value = -5555.5
$ws.Range("H97").Value = 9273326
$ws.Range("I97").Value = 1277.2
$ws.Range("K97").Value = 1277.2
$ws.Range("M97").Value = -781.2
$ws.Range("H132").Value = 6428.5776
$ws.Range("I132").Value = 4933.355
$ws.Range("K132").Value = 14800.065
$ws.Range("M132").Value = -12270.065

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6412514
$ws.Range("I20").Value = 8773491
$ws.Range("J20").Value = 4146.7144
$ws.Range("K20").Value = 8773491
$ws.Range("L20").Value = 4146.7144
$ws.Range("M20").Value = -8773244
$ws.Range("N20").Value = -4640.7144
$ws.Range("J64").Value = 2153.9375
$ws.Range("L64").Value = 2153.9375
$ws.Range("N64").Value = -2603.9375
$ws.Range("J67").Value = 2153.9375
$ws.Range("L67").Value = 2153.9375
$ws.Range("N67").Value = -3713.9375
$ws.Range("H82").Value = 2950
$ws.Range("I82").Value = 2950
$ws.Range("K82").Value = 2950
$ws.Range("M82").Value = -2567
$ws.Range("H85").Value = 2950
$ws.Range("I85").Value = 2950
$ws.Range("K85").Value = 2950
$ws.Range("M85").Value = -1624
$ws.Range("H86").Value = 43481800
$ws.Range("I86").Value = 2812.8948
$ws.Range("K86").Value = 2812.8948
$ws.Range("M86").Value = -1689.8948
$ws.Range("H89").Value = 43481800
$ws.Range("I89").Value = 2812.8948
$ws.Range("K89").Value = 14064.474
$ws.Range("M89").Value = -8448.474
$ws.Range("H94").Value = 939.375
$ws.Range("J94").Value = 2390
$ws.Range("L94").Value = 2390
$ws.Range("N94").Value = -3292
$ws.Range("H105").Value = 2533.617
$ws.Range("I105").Value = 2232.2368
$ws.Range("J105").Value = 3806.111
$ws.Range("K105").Value = 2232.2368
$ws.Range("L105").Value = 3806.111
$ws.Range("M105").Value = -485.2368000000001
$ws.Range("N105").Value = -7300.111
$ws.Range("H134").Value = 4866.18
$ws.Range("I134").Value = 3043.468
$ws.Range("J134").Value = 10985.286
$ws.Range("K134").Value = 9130.403999999999
$ws.Range("L134").Value = 32955.858
$ws.Range("M134").Value = -6595.403999999999
$ws.Range("N134").Value = -38025.858

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7070.603
$ws.Range("I31").Value = 3020.303
$ws.Range("K31").Value = 3020.303
$ws.Range("M31").Value = -2725.303
$ws.Range("H34").Value = 7070.603
$ws.Range("I34").Value = 3020.303
$ws.Range("K34").Value = 3020.303
$ws.Range("M34").Value = -2818.303
$ws.Range("H58").Value = 13164123
$ws.Range("I58").Value = 27779926
$ws.Range("K58").Value = 27779926
$ws.Range("M58").Value = -27779723
$ws.Range("H122").Value = 1498.3334
$ws.Range("I122").Value = 1269.5
$ws.Range("K122").Value = 3808.5
$ws.Range("M122").Value = -1358.5
$ws.Range("H136").Value = 13164123
$ws.Range("I136").Value = 27779926
$ws.Range("K136").Value = 83339778
$ws.Range("M136").Value = -83337228
$ws.Range("H139").Value = 52125.57
$ws.Range("J139").Value = 53313.168
$ws.Range("L139").Value = 53313.168
$ws.Range("N139").Value = -63593.168

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 3335866.2
$ws.Range("J5").Value = 4733
$ws.Range("L5").Value = 14199
$ws.Range("N5").Value = -14423
$ws.Range("H92").Value = 7693738.5
$ws.Range("J92").Value = 7693738.5
$ws.Range("L92").Value = 23081215.5
$ws.Range("N92").Value = -23083711.5
$ws.Range("H113").Value = 3632.818
$ws.Range("J113").Value = 3995.7778
$ws.Range("L113").Value = 11987.3334
$ws.Range("N113").Value = -16327.3334
$ws.Range("H121").Value = 25000778
$ws.Range("I121").Value = 20000560
$ws.Range("K121").Value = 60001680
$ws.Range("M121").Value = -60000370
$ws.Range("H129").Value = 92031.73
$ws.Range("I129").Value = 1138.8889
$ws.Range("J129").Value = 501049.5
$ws.Range("K129").Value = 3416.6667
$ws.Range("L129").Value = 1503148.5
$ws.Range("M129").Value = 1583.3333
$ws.Range("N129").Value = -1513148.5
$ws.Range("H135").Value = 3335866.2
$ws.Range("J135").Value = 4733
$ws.Range("L135").Value = 42597
$ws.Range("N135").Value = -47667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2620.96
$ws.Range("I97").Value = 2401.7856
$ws.Range("K97").Value = 2401.7856
$ws.Range("M97").Value = -1905.7856
$ws.Range("H132").Value = 2760.1428
$ws.Range("I132").Value = 2732.6667
$ws.Range("K132").Value = 8198.000100000001
$ws.Range("M132").Value = -5668.000100000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7422.5386
$ws.Range("I7").Value = 6294
$ws.Range("K7").Value = 6294
$ws.Range("M7").Value = -6182
$ws.Range("H40").Value = 5708.407
$ws.Range("I40").Value = 5643.316
$ws.Range("J40").Value = 5863
$ws.Range("K40").Value = 5643.316
$ws.Range("L40").Value = 5863
$ws.Range("M40").Value = -5507.316
$ws.Range("N40").Value = -6135
$ws.Range("H55").Value = 62500416
$ws.Range("J55").Value = 495.3846
$ws.Range("L55").Value = 495.3846
$ws.Range("N55").Value = -841.3846
$ws.Range("H126").Value = 7422.5386
$ws.Range("I126").Value = 6294
$ws.Range("K126").Value = 18882
$ws.Range("M126").Value = -16412
$ws.Range("H136").Value = 7975.548
$ws.Range("I136").Value = 4396.7827
$ws.Range("J136").Value = 12307.737
$ws.Range("K136").Value = 13190.3481
$ws.Range("L136").Value = 36923.211
$ws.Range("M136").Value = -10640.3481
$ws.Range("N136").Value = -42023.211

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2164.75
$ws.Range("I96").Value = 2059.6667
$ws.Range("K96").Value = 2059.6667
$ws.Range("M96").Value = -686.6667000000002
$ws.Range("H126").Value = 3085.6155
$ws.Range("I126").Value = 1990.7693
$ws.Range("J126").Value = 4180.4614
$ws.Range("K126").Value = 5972.3079
$ws.Range("L126").Value = 12541.3842
$ws.Range("M126").Value = -3502.3079
$ws.Range("N126").Value = -17481.3842
